$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (Price + Volume(1h)) per Jan 26 2023 GitHub Actions refresh.
# Cells store text (e.g. "306.86", "2.01%"), so force text format before writing
# to avoid Excel auto-converting the literals to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.75%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.063"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.53%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08086"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.47%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.945"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.146"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.20%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.808"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.19%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9355"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.38%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1328"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.79%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1923"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.42%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09256"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.09%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03521"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.06%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09881"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001416"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.72%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005789"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.23%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.606"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.55%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.921"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.32%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3429"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.73%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1333"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.88%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.187"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.74%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2617"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.73%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04393"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.13%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.78%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004774"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001306"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.98%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003128"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.86%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01993"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.28%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05030"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.69%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01124"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "15.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007632"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.77%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1378"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.05%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002110"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.19%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01134"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "21.45%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006394"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.04%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.41%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.15%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001190"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-28.52%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.41%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.41%"
